# Append a new data row (row 91) to the PEBCOM sheet, matching the
# incoming automated-update feed entry for case -585.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 91

function Set-TextCell($r, $c, $value) {
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $row 1 "-585"
Set-TextCell $row 2 "9/8/2025"
Set-TextCell $row 3 "Rio Cuarto 3267"
Set-TextCell $row 4 "4"
Set-TextCell $row 5 "Pendente ADM"
Set-TextCell $row 6 "PEBCOM"
Set-TextCell $row 7 "Pendiente"
Set-TextCell $row 8 "Desmonte de columna"

$ws.Cells.Item($row, 9).Value = 1

Set-TextCell $row 10 "Desmonte"
Set-TextCell $row 11 "Sin equipos"
Set-TextCell $row 12 "Terminal"

$ws.Cells.Item($row, 13).Value = -58.39368
$ws.Cells.Item($row, 14).Value = -34.652663

Set-TextCell $row 15 "San Telmo"
Set-TextCell $row 16 "Capital Sur"
